$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Rename the sheet from "repayment_20250923_20250923 (4)" to "repayment_20250923_20250923 (6)"
$ws.Name = "repayment_20250923_20250923 (6)"

# Update cell values per the diff
$ws.Range("H2").Value = 1.759

$ws.Range("D3").Value = 3
Set-TextValue "E3" "720,219.00"
Set-TextValue "G3" "0.40"
$ws.Range("H3").Value = 395

$ws.Range("H4").Value = 118

$ws.Range("D5").Value = 24
Set-TextValue "E5" "7,542,074.00"
Set-TextValue "G5" "4.53"
$ws.Range("H5").Value = 967

$ws.Range("H6").Value = 186

$ws.Range("D7").Value = 3
Set-TextValue "E7" "2,254,791.00"
Set-TextValue "G7" "1.33"
$ws.Range("H7").Value = 1.84

$ws.Range("H8").Value = 695

$ws.Range("H9").Value = 969
$ws.Range("J9").Value = 1
Set-TextValue "K9" "3.98"
Set-TextValue "L9" "5.88"

$ws.Range("H10").Value = 511

$ws.Range("D11").Value = 2
Set-TextValue "E11" "1,011,757.00"
Set-TextValue "G11" "0.72"
$ws.Range("H11").Value = 866

$ws.Range("H12").Value = 542

$ws.Range("D13").Value = 3
Set-TextValue "E13" "751,493.00"
Set-TextValue "G13" "0.54"
$ws.Range("H13").Value = 139

$ws.Range("H14").Value = 354

$ws.Range("H15").Value = 144

$ws.Range("H16").Value = 100

$ws.Range("H17").Value = 675

$ws.Range("H18").Value = 299

# Update the selection to match the new view state
$ws.Range("A2:A18").Select
